# Update cryptocurrency price (D) and volume-change (E) columns
# to the latest scraped values, per the scheduled GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.376.56"
$ws.Range("E2").Value = "  +2.30%  "
$ws.Range("D3").Value = "2.092.14"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  -1.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "342.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("E6").Value = "  -0.84%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5234"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.51%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4425"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.80%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.48"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09334"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.99%  "
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.87"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.600"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.900"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.28%  "
$ws.Range("D15").Value = "2.039.49"
$ws.Range("E15").Value = "  -2.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "101.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001159"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "21.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06665"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.331"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.27%  "
$ws.Range("E22").Value = "  -0.76%  "
$ws.Range("D23").Value = "30.342.97"
$ws.Range("E23").Value = "  +2.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.55"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.299"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.08%  "
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.42%  "
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "133.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.141"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.675"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.58%  "
$ws.Range("E32").Value = "  -0.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.761"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +9.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.248"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.863"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.19"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02634"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06837"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6986"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.343"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.52"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.61%  "
$ws.Range("E42").Value = "  -0.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6825"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.348"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.373"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.630"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000345"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.57%  "
$ws.Range("E50").Value = "  +8.04%  "
